# Update cosinor analysis results (CircaDB / CircadiPy re-run) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("E2").Value = 24.08000000000033
$ws.Range("H2").Value = [double]"1.769279720518178e-16"
$ws.Range("K2").Value = 48.6741721775787
$ws.Range("L2").Value = "[40.85606489695125, 56.492279458206156]"
$ws.Range("O2").Value = 1.62897396852804
$ws.Range("P2").Value = "[1.452868674633116, 1.8050792624229643]"
$ws.Range("S2").Value = 60.5775533551185
$ws.Range("T2").Value = "[55.511764414475365, 65.64334229576164]"
$ws.Range("W2").Value = 17.83703703703728
$ws.Range("X2").Value = 17.16212212212235
$ws.Range("Y2").Value = 18.51195195195221

# ---- Row 3 ----
$ws.Range("E3").Value = 23.42000000000022
$ws.Range("H3").Value = [double]"1.769279720518178e-16"
$ws.Range("K3").Value = 51.33129916732155
$ws.Range("L3").Value = "[40.33012729100157, 62.332471043641526]"
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0.2830263651882694
$ws.Range("P3").Value = "[0.05660527303765317, 0.5094474573388856]"
$ws.Range("Q3").Value = 0.01450859733313026
$ws.Range("R3").Value = 0.01450859733313026
$ws.Range("S3").Value = 64.35201009581334
$ws.Range("T3").Value = "[58.01248400614553, 70.69153618548116]"
$ws.Range("W3").Value = 22.36504504504526
$ws.Range("X3").Value = 21.52108108108128
$ws.Range("Y3").Value = 23.20900900900923

Write-Output "edit applied"
